$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 115.583336
$ws.Range("I33").Value = 106.7
$ws.Range("K33").Value = 106.7
$ws.Range("M33").Value = 122.3
$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("H40").Value = 1991.3334
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 1991.3334
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 1991.3334
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -2341.3334
$ws.Range("H43").Value = 926.1429000000001
$ws.Range("J43").Value = 1998
$ws.Range("L43").Value = 1998
$ws.Range("N43").Value = -2136
$ws.Range("H107").Value = 2862.7144
$ws.Range("I107").Value = 1007.8
$ws.Range("J107").Value = 7500
$ws.Range("K107").Value = 1007.8
$ws.Range("L107").Value = 7500
$ws.Range("M107").Value = 912.2
$ws.Range("N107").Value = -11340
$ws.Range("H118").Value = 403.2
$ws.Range("I118").Value = 403.2
$ws.Range("K118").Value = 1209.6
$ws.Range("M118").Value = 447.4000000000001
$ws.Range("H125").Value = 3828.6924
$ws.Range("I125").Value = 2177.5
$ws.Range("J125").Value = 9332.666999999999
$ws.Range("K125").Value = 19597.5
$ws.Range("L125").Value = 83994.003
$ws.Range("M125").Value = -17137.5
$ws.Range("N125").Value = -88914.003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4989.2964
$ws.Range("I32").Value = 4989.2964
$ws.Range("K32").Value = 4989.2964
$ws.Range("M32").Value = -4702.2964
$ws.Range("H132").Value = 1966.2222
$ws.Range("I132").Value = 1799.3846
$ws.Range("K132").Value = 5398.1538
$ws.Range("M132").Value = -2868.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2889.2593
$ws.Range("I20").Value = 2245.7058
$ws.Range("J20").Value = 3983.3
$ws.Range("K20").Value = 2245.7058
$ws.Range("L20").Value = 3983.3
$ws.Range("M20").Value = -1998.7058
$ws.Range("N20").Value = -4477.3
$ws.Range("H63").Value = 48000
$ws.Range("I63").Value = 48000
$ws.Range("K63").Value = 48000
$ws.Range("M63").Value = -47314
$ws.Range("H66").Value = 48000
$ws.Range("I66").Value = 48000
$ws.Range("K66").Value = 144000
$ws.Range("M66").Value = -140568
$ws.Range("H68").Value = 125295
$ws.Range("J68").Value = 125295
$ws.Range("L68").Value = 125295
$ws.Range("N68").Value = -126917
$ws.Range("H71").Value = 125295
$ws.Range("J71").Value = 125295
$ws.Range("L71").Value = 375885
$ws.Range("N71").Value = -383997
$ws.Range("H76").Value = 20293.2
$ws.Range("J76").Value = 20293.2
$ws.Range("L76").Value = 20293.2
$ws.Range("N76").Value = -20923.2
$ws.Range("H79").Value = 20293.2
$ws.Range("J79").Value = 20293.2
$ws.Range("L79").Value = 20293.2
$ws.Range("N79").Value = -22477.2
$ws.Range("H86").Value = 3498.3333
$ws.Range("I86").Value = 3498.3333
$ws.Range("K86").Value = 3498.3333
$ws.Range("M86").Value = -2375.3333
$ws.Range("H89").Value = 3498.3333
$ws.Range("I89").Value = 3498.3333
$ws.Range("K89").Value = 17491.6665
$ws.Range("M89").Value = -11875.6665
$ws.Range("H94").Value = 1756.5
$ws.Range("I94").Value = 1806.8572
$ws.Range("K94").Value = 1806.8572
$ws.Range("M94").Value = -1355.8572
$ws.Range("H106").Value = 7889.6665
$ws.Range("J106").Value = 7889.6665
$ws.Range("L106").Value = 7889.6665
$ws.Range("N106").Value = -10413.6665
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 5217.6665
$ws.Range("I134").Value = 5405
$ws.Range("K134").Value = 16215
$ws.Range("M134").Value = -13680

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 340.25
$ws.Range("I2").Value = 370.33334
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 370.33334
$ws.Range("L2").Value = 250
$ws.Range("M2").Value = -257.33334
$ws.Range("N2").Value = -476
$ws.Range("H7").Value = 199.5
$ws.Range("I7").Value = 199
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 199
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = -86
$ws.Range("N7").Value = -426

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1064
$ws.Range("I131").Value = 978.75
$ws.Range("J131").Value = 1405
$ws.Range("K131").Value = 2936.25
$ws.Range("L131").Value = 4215
$ws.Range("M131").Value = 2103.75
$ws.Range("N131").Value = -14295
$ws.Range("H139").Value = 2015
$ws.Range("I139").Value = 1608.6364
$ws.Range("K139").Value = 4825.9092
$ws.Range("M139").Value = 314.0907999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 132.85715
$ws.Range("I2").Value = 177
$ws.Range("J2").Value = 74
$ws.Range("K2").Value = 177
$ws.Range("L2").Value = 74
$ws.Range("M2").Value = -64
$ws.Range("N2").Value = -300
$ws.Range("H3").Value = 301
$ws.Range("I3").Value = 301
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 301
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -185
$ws.Range("N3").ClearContents()
$ws.Range("H11").Value = 6000033
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 6000033
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 6000033
$ws.Range("M11").ClearContents()
$ws.Range("N11").Value = -6000311
$ws.Range("H14").Value = 50
$ws.Range("I14").Value = 50
$ws.Range("K14").Value = 50
$ws.Range("M14").Value = 118
$ws.Range("H19").Value = 3505
$ws.Range("I19").Value = 3505
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 3505
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = -3217
$ws.Range("N19").ClearContents()
$ws.Range("H97").Value = 861.2727
$ws.Range("I97").Value = 764.55554
$ws.Range("K97").Value = 764.55554
$ws.Range("M97").Value = -268.55554
$ws.Range("H102").Value = 3665.5
$ws.Range("I102").Value = 3665.5
$ws.Range("K102").Value = 3665.5
$ws.Range("M102").Value = -2043.5
$ws.Range("H113").Value = 3680.4
$ws.Range("I113").Value = 3680.4
$ws.Range("K113").Value = 3680.4
$ws.Range("M113").Value = -1510.4
$ws.Range("H122").Value = 4801.6
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 6336
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 19008
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -23908
$ws.Range("H126").Value = 4518.8
$ws.Range("I126").Value = 4666.3335
$ws.Range("J126").Value = 4297.5
$ws.Range("K126").Value = 13999.0005
$ws.Range("L126").Value = 12892.5
$ws.Range("M126").Value = -11529.0005
$ws.Range("N126").Value = -17832.5
$ws.Range("H132").Value = 3304
$ws.Range("I132").Value = 2456
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7368
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -4838
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1866.6666
$ws.Range("I40").Value = 1866.6666
$ws.Range("K40").Value = 1866.6666
$ws.Range("M40").Value = -1730.6666
$ws.Range("H68").Value = 26000
$ws.Range("I68").Value = 4000
$ws.Range("K68").Value = 4000
$ws.Range("M68").Value = -3251
$ws.Range("H71").Value = 26000
$ws.Range("I71").Value = 4000
$ws.Range("K71").Value = 20000
$ws.Range("M71").Value = -16256
$ws.Range("H82").Value = 1447.125
$ws.Range("I82").Value = 991
$ws.Range("K82").Value = 991
$ws.Range("M82").Value = -630
$ws.Range("H85").Value = 1447.125
$ws.Range("I85").Value = 991
$ws.Range("K85").Value = 991
$ws.Range("M85").Value = 257
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()
$ws.Range("H100").Value = 1784.3334
$ws.Range("I100").Value = 1784.3334
$ws.Range("K100").Value = 1784.3334
$ws.Range("M100").Value = -1243.3334
$ws.Range("H104").Value = 21000
$ws.Range("J104").Value = 21000
$ws.Range("L104").Value = 21000
$ws.Range("N104").Value = -27988
$ws.Range("H122").Value = 3646.2856
$ws.Range("I122").Value = 3503.1667
$ws.Range("K122").Value = 10509.5001
$ws.Range("M122").Value = -8059.500100000001
$ws.Range("H132").Value = 7557.5713
$ws.Range("I132").Value = 6163
$ws.Range("K132").Value = 18489
$ws.Range("M132").Value = -15959

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 3000
$ws.Range("J19").Value = 3000
$ws.Range("L19").Value = 3000
$ws.Range("N19").Value = -3348
$ws.Range("H54").Value = 16285.714
$ws.Range("I54").Value = 10800
$ws.Range("K54").Value = 10800
$ws.Range("M54").Value = -10280
$ws.Range("H81").Value = 766.6667
$ws.Range("I81").Value = 766.6667
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1533.3334
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -472.3334
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 766.6667
$ws.Range("I84").Value = 766.6667
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 7666.666999999999
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2362.666999999999
$ws.Range("N84").ClearContents()
$ws.Range("H100").Value = 2500
$ws.Range("J100").Value = 2500
$ws.Range("L100").Value = 5000
$ws.Range("N100").Value = -6082
